# Apply the edit described by the commit: add a new row (14) of data to
# Sheet1, reusing the existing shared string "sam" (the same text already
# used in row 11) for column A, and 123 for column B. Finally move the
# active selection to N15, matching the final saved state of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A14").Value = "sam"
$ws.Range("B14").Value = 123

$ws.Range("N15").Select()
